# Insert a new record row at row 176 (pushing the existing rows 176-246 down
# to 177-247), then populate the new row with its own data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("176").Insert()

$ws.Range("A176").Value = 10
$ws.Range("B176").Value = "Vega Modelo de Temuco"
$ws.Range("C176").Value = "La Araucanía"
$ws.Range("D176").Value = 44510
$ws.Range("E176").Value = 9
$ws.Range("F176").Value = 100112037
$ws.Range("G176").Value = "Cebollín"
$ws.Range("H176").Value = "Sin especificar"
$ws.Range("I176").Value = "Primera"
$ws.Range("J176").Value = 40
$ws.Range("K176").Value = 5000
$ws.Range("L176").Value = 5000
$ws.Range("M176").Value = 5000
$ws.Range("N176").Value = "$/docena de paquetes"
$ws.Range("O176").Value = "Región de O'Higgins"
$ws.Range("P176").Value = 417
$ws.Range("Q176").Value = 12
$ws.Range("R176").Value = "Hortaliza"
